$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 72.14973449707031
$ws.Range("B3").Value = 79.62406921386719
$ws.Range("B4").Value = 79.27431488037109
$ws.Range("B5").Value = 70.43294525146484
$ws.Range("B6").Value = 73.15251159667969
$ws.Range("B7").Value = 70.60429382324219
$ws.Range("B8").Value = 73.05410003662109
$ws.Range("B9").Value = 69.86053466796875
$ws.Range("B10").Value = 74.09227752685547
$ws.Range("B11").Value = 72.92491149902344
$ws.Range("B12").Value = 75.84742736816406
$ws.Range("B13").Value = 80.71454620361328
$ws.Range("B14").Value = 86.64016723632812
$ws.Range("B15").Value = 90.40223693847656
$ws.Range("B16").Value = 111.0016326904297
$ws.Range("B17").Value = 132.1978454589844
$ws.Range("B18").Value = 126.7549514770508
$ws.Range("B19").Value = 138.9434967041016
$ws.Range("B20").Value = 137.4134216308594
$ws.Range("B21").Value = 141.7583770751953
$ws.Range("B22").Value = 146.3133239746094
$ws.Range("B23").Value = 140.6443176269531
$ws.Range("B24").Value = 148.7510375976562
$ws.Range("B25").Value = 145.7429351806641
$ws.Range("B26").Value = 144.8786163330078
$ws.Range("B27").Value = 146.5385589599609
$ws.Range("B28").Value = 136.9938049316406
$ws.Range("B29").Value = 131.2630157470703
$ws.Range("B30").Value = 130.9184112548828
$ws.Range("B31").Value = 135.4780426025391
$ws.Range("B32").Value = 137.4062194824219
$ws.Range("B33").Value = 164.2551574707031
$ws.Range("B34").Value = 153.3765258789062
$ws.Range("B35").Value = 180.8703918457031
$ws.Range("B36").Value = 171.3769989013672
$ws.Range("B37").Value = 174.994873046875
$ws.Range("B38").Value = 154.9044494628906
$ws.Range("B39").Value = 142.5528564453125
$ws.Range("B40").Value = 143.8388977050781
$ws.Range("B41").Value = 127.1858978271484
$ws.Range("B42").Value = 120.0325698852539
$ws.Range("B43").Value = 104.9885787963867
$ws.Range("B44").Value = 114.2271118164062
$ws.Range("B45").Value = 100.7412414550781
$ws.Range("B46").Value = 104.1725921630859
$ws.Range("B47").Value = 89.85700988769531
$ws.Range("B48").Value = 95.99959564208984
$ws.Range("B49").Value = 94.88642120361328
